# Update the "Förändrad" (Changed) date column (column C) for every data
# row on the active sheet from 2023-09-09 (serial 45178) to 2023-09-10
# (serial 45179).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row  # xlUp = -4162

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45178) {
        $cell.Value = 45179
    }
}
